$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume figures to match the latest scrape.
# A couple of "Price" cells are numeric-looking strings whose trailing
# zero would be silently dropped if Excel auto-coerced them to a
# number (e.g. "208.50" -> 208.5), so those are pinned to Text format
# first, exactly like Excel does when you pre-format a cell as Text.
$ws.Range("D2").Value = "26.129.42"
$ws.Range("E2").Value = "  -2.34%  "
$ws.Range("D3").Value = "1.571.88"
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.50"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("E6").Value = "  -3.60%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  -1.86%  "
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").Value = "19.55"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "1.793.42"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.594.93"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "0.512"
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("D16").Value = "64.31"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "26.112.29"
$ws.Range("E17").Value = "  -2.28%  "
$ws.Range("E18").Value = "  -1.94%  "
$ws.Range("D19").Value = "7.26"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").Value = "206.82"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("E23").Value = "  -1.50%  "
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("D25").Value = "144.08"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").Value = "6.96"
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D29").Value = "15.21"
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.20"
$ws.Range("D34").Value = "1.279.08"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("D39").Value = "1.07"
$ws.Range("E39").Value = "  -10.06%  "
$ws.Range("E41").Value = "  +2.60%  "
$ws.Range("D42").Value = "2.13"
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("E43").Value = "  -2.44%  "
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").Value = "1.705.72"
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("D46").Value = "89.01"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").Value = "0.0₆0104"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("E48").Value = "  -2.28%  "
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("E51").Value = "  -0.17%  "
